# delete akm fire selector
# Row 38 on the active worksheet is "izhmash_akm_fire_selector" / "Izhmash AKM Fire Selector".
# Deleting the entire row shifts all subsequent rows up by one, which matches the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(38).Delete()

# Update the cursor/selection to match the author's final position in the file.
$ws.Application.Goto($ws.Range("H31"))
